# Auto-generated edit script: updates market-data driven profit columns
# (H..N) across all 8 job sheets, matching the scheduled runner's refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

$updates = @(
    @("I33", 17222.895),
    @("J33", 2024.6666),
    @("K33", 17222.895),
    @("L33", 2024.6666),
    @("M33", -16993.895),
    @("N33", -2482.6666),
    @("H51", 8708.333000000001),
    @("J51", 8818.182000000001),
    @("L51", 8818.182000000001),
    @("N51", -9786.182000000001),
    @("H80", 275.11765),
    @("I80", 294.77777),
    @("J80", 253),
    @("K80", 884.33331),
    @("L80", 759),
    @("M80", 113.66669),
    @("N80", -2755),
    @("H83", 275.11765),
    @("I83", 294.77777),
    @("J83", 253),
    @("K83", 2652.99993),
    @("L83", 2277),
    @("M83", 2339.00007),
    @("N83", -12261),
    @("H132", 2074.4119),
    @("I132", 1628.0769),
    @("K132", 4884.2307),
    @("M132", -2354.2307),
    @("H138", 26317776),
    @("I138", 1916.5294),
    @("J138", 47621090),
    @("K138", 5749.5882),
    @("L138", 142863270),
    @("M138", -609.5882000000001),
    @("N138", -142873550)
)
foreach ($u in $updates) {
    $ws.Range($u[0]).Value = $u[1]
}

$ws = $wb.Worksheets.Item("ARM")

$updates = @(
    @("H32", 3442.1428),
    @("I32", 3355.8552),
    @("K32", 3355.8552),
    @("M32", -3068.8552),
    @("H45", 4901.9375),
    @("I45", 2291.7144),
    @("J45", 6932.1113),
    @("K45", 2291.7144),
    @("L45", 6932.1113),
    @("M45", -1914.7144),
    @("N45", -7686.1113),
    @("H61", 50019900),
    @("I61", 100026000),
    @("K61", 100026000),
    @("M61", -100025788),
    @("H74", 2791.7666),
    @("I74", 2312.6072),
    @("K74", 2312.6072),
    @("M74", -1438.6072),
    @("H77", 2791.7666),
    @("I77", 2312.6072),
    @("K77", 11563.036),
    @("M77", -7195.036),
    @("H136", 50019900),
    @("I136", 100026000),
    @("K136", 300078000),
    @("M136", -300075450)
)
foreach ($u in $updates) {
    $ws.Range($u[0]).Value = $u[1]
}

$ws = $wb.Worksheets.Item("CRP")

$updates = @(
    @("H4", 145141.42),
    @("I4", 202198),
    @("K4", 202198),
    @("M4", -202086),
    @("H31", 5207.9375),
    @("I31", 3893.3572),
    @("K31", 3893.3572),
    @("M31", -3598.3572),
    @("H34", 5207.9375),
    @("I34", 3893.3572),
    @("K34", 3893.3572),
    @("M34", -3691.3572),
    @("H122", 3364.0952),
    @("I122", 2724.7778),
    @("K122", 8174.3334),
    @("M122", -5724.3334),
    @("H134", 5772.2593),
    @("I134", 4700.625),
    @("J134", 7331),
    @("K134", 14101.875),
    @("L134", 21993),
    @("M134", -11566.875),
    @("N134", -27063),
    @("H135", 0),
    @("I135", 0),
    @("J135", 0),
    @("K135", 0),
    @("L135", 0)
)
foreach ($u in $updates) {
    $ws.Range($u[0]).Value = $u[1]
}
$clears = @("M135", "N135")
foreach ($c in $clears) {
    $ws.Range($c).ClearContents()
}

$ws = $wb.Worksheets.Item("CUL")

$updates = @(
    @("H3", 11416.143),
    @("I3", 9075.091),
    @("K3", 27225.273),
    @("M3", -27113.273),
    @("H5", 1596.8235),
    @("I5", 1472.8462),
    @("K5", 4418.5386),
    @("M5", -4306.5386),
    @("H14", 133.1),
    @("I14", 133.1),
    @("K14", 399.3),
    @("M14", -226.3),
    @("H51", 1212.8334),
    @("I51", 996.4),
    @("K51", 2989.2),
    @("M51", -2529.2),
    @("H135", 1596.8235),
    @("I135", 1472.8462),
    @("K135", 13255.6158),
    @("M135", -10720.6158),
    @("H137", 1066),
    @("I137", 1066),
    @("K137", 3198),
    @("M137", 1902),
    @("H140", 1851.5),
    @("I140", 1447.3636),
    @("K140", 4342.0908),
    @("M140", 837.9092000000001)
)
foreach ($u in $updates) {
    $ws.Range($u[0]).Value = $u[1]
}

$ws = $wb.Worksheets.Item("GSM")

$updates = @(
    @("H80", 4922.857),
    @("I80", 5653.3335),
    @("J80", 4375),
    @("K80", 5653.3335),
    @("L80", 4375),
    @("M80", -4655.3335),
    @("N80", -6371),
    @("H83", 4922.857),
    @("I83", 5653.3335),
    @("J83", 4375),
    @("K83", 28266.6675),
    @("L83", 21875),
    @("M83", -23274.6675),
    @("N83", -31859),
    @("H93", 40666.668),
    @("J93", 40666.668),
    @("L93", 40666.668),
    @("N93", -44410.668),
    @("H97", 4354),
    @("I97", 846.5454999999999),
    @("K97", 846.5454999999999),
    @("M97", -350.5454999999999),
    @("H131", 50000),
    @("J131", 50000),
    @("L131", 50000),
    @("N131", -60080),
    @("H132", 4774.143),
    @("I132", 3944.8333),
    @("K132", 11834.4999),
    @("M132", -9304.499899999999)
)
foreach ($u in $updates) {
    $ws.Range($u[0]).Value = $u[1]
}

$ws = $wb.Worksheets.Item("LTW")

$updates = @(
    @("H82", 365.75),
    @("I82", 188),
    @("K82", 188),
    @("M82", 173),
    @("H85", 365.75),
    @("I85", 188),
    @("K85", 188),
    @("M85", 1060),
    @("H107", 3212.6667),
    @("I107", 3212.6667),
    @("K107", 3212.6667),
    @("M107", -1292.6667),
    @("H122", 2910.7778),
    @("I122", 2837.0938),
    @("K122", 8511.2814),
    @("M122", -6061.2814)
)
foreach ($u in $updates) {
    $ws.Range($u[0]).Value = $u[1]
}

$ws = $wb.Worksheets.Item("WVR")

$updates = @(
    @("H64", 0),
    @("J64", 0),
    @("L64", 0),
    @("H67", 0),
    @("J67", 0),
    @("L67", 0),
    @("H81", 4847.7),
    @("J81", 7000.6665),
    @("L81", 14001.333),
    @("N81", -16123.333),
    @("H84", 4847.7),
    @("J84", 7000.6665),
    @("L84", 70006.66500000001),
    @("N84", -80614.66500000001),
    @("H122", 3000.2144),
    @("I122", 1666.9166),
    @("K122", 5000.7498),
    @("M122", -2550.7498),
    @("H132", 2962.7112),
    @("I132", 2507.6904),
    @("K132", 7523.0712),
    @("M132", -4993.0712)
)
foreach ($u in $updates) {
    $ws.Range($u[0]).Value = $u[1]
}
$clears = @("N64", "N67")
foreach ($c in $clears) {
    $ws.Range($c).ClearContents()
}
